$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new status entry for erica (hours update).
# Column A holds plain-text dates (no special number format, like the
# existing rows above it). Entering "2/6/2010" directly would make Excel
# auto-convert it to a date serial, so instead build it as a text formula
# and then paste-special just the value back over itself; that commits the
# literal text to the cell without leaving any date number format behind.
$ws.Range("A20").Formula = "=""2/6/2010"""
$ws.Range("A20").Copy()
$ws.Range("A20").PasteSpecial(-4163)

$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "CFP Update and misc"

# Move the active selection to A21, matching the post-edit state in Excel
$ws.Range("A21").Select()
